$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Azure Networking")

# Row 19: wp-dev-repsrc
$ws.Cells.Item(19, 1).Value = "Whitespace Development and Test"
$ws.Cells.Item(19, 2).Value = "wp-dev-repsrc"
$ws.Cells.Item(19, 3).Value = "10.41.0.0/23"
$ws.Cells.Item(19, 4).Value = "app-sub 10.41.0.0/24`nappgw-sub 10.41.1.0/24"
$ws.Cells.Item(19, 6).Value = "wp-dev-repsrc-uksouth.internal"
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

# Row 20: wp-dev-repdst
$ws.Cells.Item(20, 1).Value = "Whitespace Development and Test"
$ws.Cells.Item(20, 2).Value = "wp-dev-repdst"
$ws.Cells.Item(20, 3).Value = "10.42.0.0/23"
$ws.Cells.Item(20, 4).Value = "app-sub 10.42.0.0/24`nappgw-sub 10.42.1.0/24"
$ws.Cells.Item(20, 6).Value = "wp-dev-repdst-uksouth.internal"
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(18).RowHeight

# Update the view: scroll to A13, select F12
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("F12").Select()
